$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4,8).Value2 = 169
$ws.Cells.Item(4,9).Value2 = 169
$ws.Cells.Item(4,10).Value2 = 0
$ws.Cells.Item(4,11).Value2 = 169
$ws.Cells.Item(4,12).Value2 = 0
$ws.Cells.Item(4,13).Value2 = -55
$ws.Cells.Item(4,14).ClearContents()

$ws.Cells.Item(29,8).Value2 = 5320.2
$ws.Cells.Item(29,10).Value2 = 5500
$ws.Cells.Item(29,12).Value2 = 16500
$ws.Cells.Item(29,14).Value2 = -17062

$ws.Cells.Item(39,8).Value2 = 164.875
$ws.Cells.Item(39,9).Value2 = 121.375
$ws.Cells.Item(39,10).Value2 = 251.875
$ws.Cells.Item(39,11).Value2 = 364.125
$ws.Cells.Item(39,12).Value2 = 755.625
$ws.Cells.Item(39,13).Value2 = -68.125
$ws.Cells.Item(39,14).Value2 = -1347.625

$ws.Cells.Item(87,8).Value2 = 95338.25
$ws.Cells.Item(87,10).Value2 = 95338.25
$ws.Cells.Item(87,12).Value2 = 95338.25
$ws.Cells.Item(87,14).Value2 = -97834.25

$ws.Cells.Item(90,8).Value2 = 95338.25
$ws.Cells.Item(90,10).Value2 = 95338.25
$ws.Cells.Item(90,12).Value2 = 286014.75
$ws.Cells.Item(90,14).Value2 = -298494.75

$ws.Cells.Item(100,8).Value2 = 1184.4445
$ws.Cells.Item(100,9).Value2 = 1433.6666
$ws.Cells.Item(100,11).Value2 = 1433.6666
$ws.Cells.Item(100,13).Value2 = -892.6666

$ws.Cells.Item(101,8).Value2 = 419.75
$ws.Cells.Item(101,9).Value2 = 326.33334
$ws.Cells.Item(101,11).Value2 = 979.0000200000001
$ws.Cells.Item(101,13).Value2 = 642.9999799999999

$ws.Cells.Item(112,8).Value2 = 2100
$ws.Cells.Item(112,10).Value2 = 1000
$ws.Cells.Item(112,12).Value2 = 3000
$ws.Cells.Item(112,14).Value2 = -5216

$ws.Cells.Item(137,8).Value2 = 2999.6667
$ws.Cells.Item(137,9).Value2 = 1997.5
$ws.Cells.Item(137,10).Value2 = 3153.8462
$ws.Cells.Item(137,11).Value2 = 5992.5
$ws.Cells.Item(137,12).Value2 = 9461.5386
$ws.Cells.Item(137,13).Value2 = -3442.5
$ws.Cells.Item(137,14).Value2 = -14561.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value2 = 34984
$ws.Cells.Item(74,9).Value2 = 34984
$ws.Cells.Item(74,11).Value2 = 34984
$ws.Cells.Item(74,13).Value2 = -34110

$ws.Cells.Item(77,8).Value2 = 34984
$ws.Cells.Item(77,9).Value2 = 34984
$ws.Cells.Item(77,11).Value2 = 174920
$ws.Cells.Item(77,13).Value2 = -170552

$ws.Cells.Item(102,8).Value2 = 3597.1
$ws.Cells.Item(102,9).Value2 = 1638.2858
$ws.Cells.Item(102,11).Value2 = 1638.2858
$ws.Cells.Item(102,13).Value2 = -16.28580000000011

$ws.Cells.Item(110,8).Value2 = 750
$ws.Cells.Item(110,9).Value2 = 750
$ws.Cells.Item(110,11).Value2 = 750
$ws.Cells.Item(110,13).Value2 = 1295

$ws.Cells.Item(122,8).Value2 = 1810.6428
$ws.Cells.Item(122,9).Value2 = 1825.3334
$ws.Cells.Item(122,11).Value2 = 5476.0002
$ws.Cells.Item(122,13).Value2 = -3026.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value2 = 3596.3076
$ws.Cells.Item(86,9).Value2 = 2231.2666
$ws.Cells.Item(86,10).Value2 = 5457.727
$ws.Cells.Item(86,11).Value2 = 2231.2666
$ws.Cells.Item(86,12).Value2 = 5457.727
$ws.Cells.Item(86,13).Value2 = -1108.2666
$ws.Cells.Item(86,14).Value2 = -7703.727

$ws.Cells.Item(89,8).Value2 = 3596.3076
$ws.Cells.Item(89,9).Value2 = 2231.2666
$ws.Cells.Item(89,10).Value2 = 5457.727
$ws.Cells.Item(89,11).Value2 = 11156.333
$ws.Cells.Item(89,12).Value2 = 27288.635
$ws.Cells.Item(89,13).Value2 = -5540.332999999999
$ws.Cells.Item(89,14).Value2 = -38520.63499999999

$ws.Cells.Item(99,8).Value2 = 2482.4443
$ws.Cells.Item(99,9).Value2 = 1757.3334
$ws.Cells.Item(99,11).Value2 = 1757.3334
$ws.Cells.Item(99,13).Value2 = -259.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value2 = 93405.45
$ws.Cells.Item(16,9).Value2 = 101946
$ws.Cells.Item(16,11).Value2 = 101946
$ws.Cells.Item(16,13).Value2 = -101659

$ws.Cells.Item(113,8).Value2 = 93405.45
$ws.Cells.Item(113,9).Value2 = 101946
$ws.Cells.Item(113,11).Value2 = 101946
$ws.Cells.Item(113,13).Value2 = -99776

$ws.Cells.Item(122,8).Value2 = 2067.7693
$ws.Cells.Item(122,9).Value2 = 1650.3
$ws.Cells.Item(122,10).Value2 = 3459.3333
$ws.Cells.Item(122,11).Value2 = 4950.9
$ws.Cells.Item(122,12).Value2 = 10377.9999
$ws.Cells.Item(122,13).Value2 = -2500.9
$ws.Cells.Item(122,14).Value2 = -15277.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2,8).Value2 = 119.2
$ws.Cells.Item(2,9).Value2 = 131.42857
$ws.Cells.Item(2,11).Value2 = 788.57142
$ws.Cells.Item(2,13).Value2 = -675.57142

$ws.Cells.Item(7,8).Value2 = 295
$ws.Cells.Item(7,9).Value2 = 100
$ws.Cells.Item(7,10).Value2 = 490
$ws.Cells.Item(7,11).Value2 = 300
$ws.Cells.Item(7,12).Value2 = 1470
$ws.Cells.Item(7,13).Value2 = -188
$ws.Cells.Item(7,14).Value2 = -1694

$ws.Cells.Item(34,8).Value2 = 1813
$ws.Cells.Item(34,10).Value2 = 725
$ws.Cells.Item(34,12).Value2 = 2175
$ws.Cells.Item(34,14).Value2 = -2343

$ws.Cells.Item(39,8).Value2 = 5179
$ws.Cells.Item(39,10).Value2 = 5179
$ws.Cells.Item(39,12).Value2 = 15537
$ws.Cells.Item(39,14).Value2 = -16125

$ws.Cells.Item(52,8).Value2 = 1162.5
$ws.Cells.Item(52,10).Value2 = 1162.5
$ws.Cells.Item(52,12).Value2 = 3487.5
$ws.Cells.Item(52,14).Value2 = -4019.5

$ws.Cells.Item(55,8).Value2 = 4638.6665
$ws.Cells.Item(55,10).Value2 = 4638.6665
$ws.Cells.Item(55,12).Value2 = 13915.9995
$ws.Cells.Item(55,14).Value2 = -14269.9995

$ws.Cells.Item(117,8).Value2 = 601
$ws.Cells.Item(117,10).Value2 = 737.7143
$ws.Cells.Item(117,12).Value2 = 2213.1429
$ws.Cells.Item(117,14).Value2 = -9097.142899999999

$ws.Cells.Item(118,8).Value2 = 0
$ws.Cells.Item(118,9).Value2 = 0
$ws.Cells.Item(118,11).Value2 = 0
$ws.Cells.Item(118,13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80,8).Value2 = 3814.5715
$ws.Cells.Item(80,9).Value2 = 3819
$ws.Cells.Item(80,10).Value2 = 3803.5
$ws.Cells.Item(80,11).Value2 = 3819
$ws.Cells.Item(80,12).Value2 = 3803.5
$ws.Cells.Item(80,13).Value2 = -2821
$ws.Cells.Item(80,14).Value2 = -5799.5

$ws.Cells.Item(83,8).Value2 = 3814.5715
$ws.Cells.Item(83,9).Value2 = 3819
$ws.Cells.Item(83,10).Value2 = 3803.5
$ws.Cells.Item(83,11).Value2 = 19095
$ws.Cells.Item(83,12).Value2 = 19017.5
$ws.Cells.Item(83,13).Value2 = -14103
$ws.Cells.Item(83,14).Value2 = -29001.5

$ws.Cells.Item(102,8).Value2 = 2567.6924
$ws.Cells.Item(102,9).Value2 = 1887
$ws.Cells.Item(102,11).Value2 = 1887
$ws.Cells.Item(102,13).Value2 = -265

$ws.Cells.Item(122,8).Value2 = 5147.6
$ws.Cells.Item(122,9).Value2 = 3934.5
$ws.Cells.Item(122,11).Value2 = 11803.5
$ws.Cells.Item(122,13).Value2 = -9353.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value2 = 6335.2
$ws.Cells.Item(7,9).Value2 = 4895.6665
$ws.Cells.Item(7,11).Value2 = 4895.6665
$ws.Cells.Item(7,13).Value2 = -4783.6665

$ws.Cells.Item(36,8).Value2 = 0
$ws.Cells.Item(36,10).Value2 = 0
$ws.Cells.Item(36,12).Value2 = 0
$ws.Cells.Item(36,14).ClearContents()

$ws.Cells.Item(53,8).Value2 = 0
$ws.Cells.Item(53,10).Value2 = 0
$ws.Cells.Item(53,12).Value2 = 0
$ws.Cells.Item(53,14).ClearContents()

$ws.Cells.Item(126,8).Value2 = 6335.2
$ws.Cells.Item(126,9).Value2 = 4895.6665
$ws.Cells.Item(126,11).Value2 = 14686.9995
$ws.Cells.Item(126,13).Value2 = -12216.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5,8).Value2 = 12144286
$ws.Cells.Item(5,9).Value2 = 13001000
$ws.Cells.Item(5,11).Value2 = 13001000
$ws.Cells.Item(5,13).Value2 = -13000888

$ws.Cells.Item(54,8).Value2 = 48888
$ws.Cells.Item(54,10).Value2 = 48888
$ws.Cells.Item(54,12).Value2 = 48888
$ws.Cells.Item(54,14).Value2 = -49928

$ws.Cells.Item(81,8).Value2 = 556.5714
$ws.Cells.Item(81,9).Value2 = 459.2
$ws.Cells.Item(81,10).Value2 = 800
$ws.Cells.Item(81,11).Value2 = 918.4
$ws.Cells.Item(81,12).Value2 = 1600
$ws.Cells.Item(81,13).Value2 = 142.6
$ws.Cells.Item(81,14).Value2 = -3722

$ws.Cells.Item(84,8).Value2 = 556.5714
$ws.Cells.Item(84,9).Value2 = 459.2
$ws.Cells.Item(84,10).Value2 = 800
$ws.Cells.Item(84,11).Value2 = 4592
$ws.Cells.Item(84,12).Value2 = 8000
$ws.Cells.Item(84,13).Value2 = 712
$ws.Cells.Item(84,14).Value2 = -18608

$ws.Cells.Item(107,8).Value2 = 505.25
$ws.Cells.Item(107,9).Value2 = 505.25
$ws.Cells.Item(107,11).Value2 = 1515.75
$ws.Cells.Item(107,13).Value2 = 404.25

$ws.Cells.Item(113,8).Value2 = 696.1667
$ws.Cells.Item(113,10).Value2 = 1124.5
$ws.Cells.Item(113,12).Value2 = 3373.5
$ws.Cells.Item(113,14).Value2 = -7713.5

$ws.Cells.Item(122,8).Value2 = 5213.4287
$ws.Cells.Item(122,10).Value2 = 6298
$ws.Cells.Item(122,12).Value2 = 18894
$ws.Cells.Item(122,14).Value2 = -23794

$ws.Cells.Item(126,8).Value2 = 3773.074
$ws.Cells.Item(126,9).Value2 = 2385.4736
$ws.Cells.Item(126,11).Value2 = 7156.4208
$ws.Cells.Item(126,13).Value2 = -4686.4208

$ws.Cells.Item(132,8).Value2 = 1317.8889
$ws.Cells.Item(132,9).Value2 = 1151.3077
$ws.Cells.Item(132,11).Value2 = 3453.9231
$ws.Cells.Item(132,13).Value2 = -923.9231
